# Applies the "chore: update Sheets via scheduled runner" commit.
# Updates computed market-price / profit columns (H:N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables with refreshed
# Universalis price snapshots. Values are plain numeric literals
# (no formulas in these sheets), so each changed cell is written
# directly; a few cells whose source data is no longer available
# are cleared instead (matches cells removed from the XML).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 1048.8235
$ws.Range("I2").Value = 329.25
$ws.Range("J2").Value = 1688.4445
$ws.Range("K2").Value = 329.25
$ws.Range("L2").Value = 1688.4445
$ws.Range("M2").Value = -216.25
$ws.Range("N2").Value = -1914.4445
$ws.Range("H5").Value = 1851.6364
$ws.Range("I5").Value = 1851.6364
$ws.Range("K5").Value = 1851.6364
$ws.Range("M5").Value = -1736.6364
$ws.Range("H9").Value = 107.92308
$ws.Range("J9").Value = 113.333336
$ws.Range("L9").Value = 113.333336
$ws.Range("N9").Value = -451.333336
$ws.Range("H12").Value = 327.7
$ws.Range("I12").Value = 331.8889
$ws.Range("J12").Value = 290
$ws.Range("K12").Value = 331.8889
$ws.Range("L12").Value = 290
$ws.Range("M12").Value = -161.8889
$ws.Range("N12").Value = -630
$ws.Range("H15").Value = 860.4643
$ws.Range("I15").Value = 860.4643
$ws.Range("K15").Value = 2581.3929
$ws.Range("M15").Value = -2412.3929
$ws.Range("H29").Value = 363.2857
$ws.Range("I29").Value = 178.6
$ws.Range("J29").Value = 825
$ws.Range("K29").Value = 535.8
$ws.Range("L29").Value = 2475
$ws.Range("M29").Value = -254.8
$ws.Range("N29").Value = -3037
$ws.Range("H33").Value = 455
$ws.Range("I33").Value = 178.33333
$ws.Range("K33").Value = 178.33333
$ws.Range("M33").Value = 50.66667000000001
$ws.Range("H40").Value = 3606.8
$ws.Range("J40").Value = 3606.8
$ws.Range("L40").Value = 3606.8
$ws.Range("N40").Value = -3956.8
$ws.Range("H100").Value = 4858.1665
$ws.Range("I100").Value = 2042.5714
$ws.Range("J100").Value = 8800
$ws.Range("K100").Value = 2042.5714
$ws.Range("L100").Value = 8800
$ws.Range("M100").Value = -1501.5714
$ws.Range("N100").Value = -9882
$ws.Range("H103").Value = 811.1429000000001
$ws.Range("J103").Value = 879.6667
$ws.Range("L103").Value = 2639.0001
$ws.Range("N103").Value = -3811.0001
$ws.Range("H111").Value = 1458.1666
$ws.Range("I111").Value = 1109.8
$ws.Range("K111").Value = 3329.4
$ws.Range("M111").Value = -262.3999999999996
$ws.Range("H137").Value = 2862.5
$ws.Range("I137").Value = 2372
$ws.Range("J137").Value = 3353
$ws.Range("K137").Value = 7116
$ws.Range("L137").Value = 10059
$ws.Range("M137").Value = -4566
$ws.Range("N137").Value = -15159

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1671.1428
$ws.Range("I2").Value = 1449.6666
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1449.6666
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -1336.6666
$ws.Range("N2").Value = -3226
$ws.Range("H32").Value = 5534.976
$ws.Range("I32").Value = 4450.4634
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 4450.4634
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -4163.4634
$ws.Range("N32").Value = -50574
$ws.Range("H45").Value = 3386.0908
$ws.Range("I45").Value = 3407.8333
$ws.Range("K45").Value = 3407.8333
$ws.Range("M45").Value = -3030.8333
$ws.Range("H61").Value = 1865.4
$ws.Range("I61").Value = 1831.75
$ws.Range("K61").Value = 1831.75
$ws.Range("M61").Value = -1619.75
$ws.Range("H74").Value = 1609.5333
$ws.Range("J74").Value = 2300
$ws.Range("L74").Value = 2300
$ws.Range("N74").Value = -4048
$ws.Range("H77").Value = 1609.5333
$ws.Range("J77").Value = 2300
$ws.Range("L77").Value = 11500
$ws.Range("N77").Value = -20236
$ws.Range("H97").Value = 1102.125
$ws.Range("I97").Value = 1136.3334
$ws.Range("J97").Value = 999.5
$ws.Range("K97").Value = 1136.3334
$ws.Range("L97").Value = 999.5
$ws.Range("M97").Value = -640.3334
$ws.Range("N97").Value = -1991.5
$ws.Range("H116").Value = 1671.1428
$ws.Range("I116").Value = 1449.6666
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 1449.6666
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 844.3334
$ws.Range("N116").Value = -7588
$ws.Range("H132").Value = 4055.318
$ws.Range("I132").Value = 3986
$ws.Range("J132").Value = 4176.625
$ws.Range("K132").Value = 11958
$ws.Range("L132").Value = 12529.875
$ws.Range("M132").Value = -9428
$ws.Range("N132").Value = -17589.875
$ws.Range("H136").Value = 1865.4
$ws.Range("I136").Value = 1831.75
$ws.Range("K136").Value = 5495.25
$ws.Range("M136").Value = -2945.25
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1671.1428
$ws.Range("I3").Value = 1449.6666
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 1449.6666
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -1335.6666
$ws.Range("N3").Value = -3228
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H94").Value = 1785.2916
$ws.Range("I94").Value = 1871.85
$ws.Range("K94").Value = 1871.85
$ws.Range("M94").Value = -1420.85
$ws.Range("H99").Value = 4109.727
$ws.Range("J99").Value = 2666.3333
$ws.Range("L99").Value = 2666.3333
$ws.Range("N99").Value = -5662.3333

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3445.4375
$ws.Range("I31").Value = 2414.2
$ws.Range("J31").Value = 5164.1665
$ws.Range("K31").Value = 2414.2
$ws.Range("L31").Value = 5164.1665
$ws.Range("M31").Value = -2119.2
$ws.Range("N31").Value = -5754.1665
$ws.Range("H34").Value = 3445.4375
$ws.Range("I34").Value = 2414.2
$ws.Range("J34").Value = 5164.1665
$ws.Range("K34").Value = 2414.2
$ws.Range("L34").Value = 5164.1665
$ws.Range("M34").Value = -2212.2
$ws.Range("N34").Value = -5568.1665
$ws.Range("H50").Value = 29996.875
$ws.Range("J50").Value = 29996.875
$ws.Range("L50").Value = 29996.875
$ws.Range("N50").Value = -31246.875
$ws.Range("H58").Value = 4171.4
$ws.Range("I58").Value = 4457.222
$ws.Range("K58").Value = 4457.222
$ws.Range("M58").Value = -4254.222
$ws.Range("H74").Value = 39998.332
$ws.Range("J74").Value = 39998.332
$ws.Range("L74").Value = 39998.332
$ws.Range("N74").Value = -41746.332
$ws.Range("H77").Value = 39998.332
$ws.Range("J77").Value = 39998.332
$ws.Range("L77").Value = 119994.996
$ws.Range("N77").Value = -128730.996
$ws.Range("H105").Value = 1273.5555
$ws.Range("I105").Value = 1282.875
$ws.Range("K105").Value = 1282.875
$ws.Range("M105").Value = 464.125
$ws.Range("H132").Value = 2833.3333
$ws.Range("I132").Value = 3250
$ws.Range("K132").Value = 9750
$ws.Range("M132").Value = -7220
$ws.Range("H136").Value = 4171.4
$ws.Range("I136").Value = 4457.222
$ws.Range("K136").Value = 13371.666
$ws.Range("M136").Value = -10821.666

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H64").Value = 267.5
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 267.5
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 1547.3572
$ws.Range("J122").Value = 1560.875
$ws.Range("L122").Value = 14047.875
$ws.Range("N122").Value = -18947.875

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H9").Value = 2507.4443
$ws.Range("J9").Value = 10032.5
$ws.Range("L9").Value = 10032.5
$ws.Range("N9").Value = -10372.5
$ws.Range("H132").Value = 741.9286
$ws.Range("I132").Value = 741.9286
$ws.Range("K132").Value = 2225.7858
$ws.Range("M132").Value = 304.2142000000003

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H9").Value = 627.6
$ws.Range("I9").Value = 362.66666
$ws.Range("J9").Value = 1025
$ws.Range("K9").Value = 362.66666
$ws.Range("L9").Value = 1025
$ws.Range("M9").Value = -138.66666
$ws.Range("N9").Value = -1473
$ws.Range("H46").Value = 4618.8
$ws.Range("I46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("M46").Value = -4812
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 2997.3333
$ws.Range("I100").Value = 2997.3333
$ws.Range("K100").Value = 2997.3333
$ws.Range("M100").Value = -2456.3333
$ws.Range("H132").Value = 24800.8
$ws.Range("I132").Value = 29001.334
$ws.Range("K132").Value = 87004.00199999999
$ws.Range("M132").Value = -84474.00199999999
$ws.Range("H136").Value = 3762.111
$ws.Range("I136").Value = 3430.5715
$ws.Range("J136").Value = 3973.0908
$ws.Range("K136").Value = 10291.7145
$ws.Range("L136").Value = 11919.2724
$ws.Range("M136").Value = -7741.7145
$ws.Range("N136").Value = -17019.2724

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H136").Value = 2661.25
$ws.Range("I136").Value = 2705.3333
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 8115.999899999999
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -5565.999899999999
$ws.Range("N136").Value = -11100

